$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run of averaged-intensity generation: Gaussian-Quadrature moved up (now right
# after "Ring Perpendicular to TD"), three new spiral sampling schemes inserted
# after it, and the three HexGrid rows that got pushed off the end are re-added
# at the bottom of the table.

$labels = @(
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "Gaussian-Quadrature",
    "Spiral-90deg-10rot-5space",
    "Spiral-90deg-15rot-5space",
    "Spiral-90deg-10rot-3space",
    "NoRotation-tilt60deg",
    "Rotation-NoTilt",
    "Rotation-60detTilt",
    "HexGrid-90degTilt5degRes",
    "HexGrid-90degTilt22p5degRes",
    "HexGrid-60degTilt5degRes"
)

# Copy the formatting of the last existing data row down onto the three brand
# new rows before filling in their values.
$ws.Range("A16:P16").Copy() | Out-Null
$ws.Range("A17:P19").PasteSpecial(-4122) | Out-Null

$startRow = 3
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $labels[$i]
    for ($col = 3; $col -le 16; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}
